$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.286.02'
$ws.Range("E2").Value = '  -0.60%  '

$ws.Range("D3").Value = '3.504.61'
$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '583.66'
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("D6").Value = '134.97'
$ws.Range("E6").Value = '  +1.88%  '

$ws.Range("D7").Value = '3.506.66'
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  -0.32%  '

$ws.Range("E10").Value = '  +0.14%  '

$ws.Range("D11").Value = '7.13'
$ws.Range("E11").Value = '  +0.24%  '

$ws.Range("E12").Value = '  -3.55%  '

$ws.Range("D13").Value = '4.101.77'
$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("E14").Value = '  -0.16%  '

$ws.Range("E15").Value = '  +1.03%  '

$ws.Range("D16").Value = '3.504.43'
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("D17").Value = '26.34'
$ws.Range("E17").Value = '  -5.17%  '

$ws.Range("D18").Value = '64.289.52'
$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("E19").Value = '  -2.65%  '

$ws.Range("D20").Value = '13.83'
$ws.Range("E20").Value = '  -2.84%  '

$ws.Range("E21").Value = '  -2.15%  '

$ws.Range("D22").Value = '383.76'
$ws.Range("E22").Value = '  -1.80%  '

$ws.Range("E23").Value = '  -1.46%  '

$ws.Range("D24").Value = '3.643.91'

$ws.Range("D25").Value = '73.83'
$ws.Range("E25").Value = '  -0.45%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").Value = '5.72'
$ws.Range("E27").Value = '  +0.22%  '

$ws.Range("E28").Value = '  +4.69%  '

$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").Value = '1.57'
$ws.Range("E29").Value = '  -0.54%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '7.56'
$ws.Range("E30").Value = '  +1.50%  '

$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").Value = '8.29'
$ws.Range("E32").Value = '  +1.17%  '

$ws.Range("E33").Value = '  -2.12%  '

$ws.Range("D34").Value = '3.523.96'
$ws.Range("E34").Value = '  +0.43%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  -0.37%  '

$ws.Range("D37").Value = '23.57'
$ws.Range("E37").Value = '  -1.59%  '

$ws.Range("E38").Value = '  +1.90%  '

$ws.Range("E39").Value = '  -3.35%  '

$ws.Range("E40").Value = '  -1.69%  '

$ws.Range("D41").Value = '164.58'
$ws.Range("E41").Value = '  -4.02%  '

$ws.Range("E42").Value = '  -4.15%  '

$ws.Range("E43").Value = '  -0.72%  '

$ws.Range("D44").Value = '25.77'
$ws.Range("E44").Value = '  -1.60%  '

$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").Value = '41.88'
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("E47").Value = '  -1.02%  '

$ws.Range("E48").Value = '  +0.07%  '

$ws.Range("E49").Value = '  -1.78%  '

$ws.Range("D50").Value = '2.471.42'
$ws.Range("E50").Value = '  -0.21%  '

$ws.Range("E51").Value = '  +1.79%  '
